# Auto-generated edit script applying F-column (想去人数) value updates
# across all 4 worksheets, per the commit diff.

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 14196
$ws.Range("F4").Value = 14196
$ws.Range("F5").Value = 14299
$ws.Range("F7").Value = 1399
$ws.Range("F9").Value = 984
$ws.Range("F15").Value = 436
$ws.Range("F16").Value = 2131
$ws.Range("F17").Value = 1198
$ws.Range("F18").Value = 1828
$ws.Range("F21").Value = 2277
$ws.Range("F22").Value = 563
$ws.Range("F23").Value = 812
$ws.Range("F24").Value = 3318
$ws.Range("F26").Value = 311
$ws.Range("F27").Value = 2391
$ws.Range("F28").Value = 592
$ws.Range("F31").Value = 1788
$ws.Range("F32").Value = 1069
$ws.Range("F33").Value = 1389
$ws.Range("F34").Value = 100
$ws.Range("F35").Value = 147
$ws.Range("F36").Value = 4810
$ws.Range("F37").Value = 4845
$ws.Range("F38").Value = 302
$ws.Range("F40").Value = 672
$ws.Range("F42").Value = 3287
$ws.Range("F45").Value = 337
$ws.Range("F46").Value = 106
$ws.Range("F47").Value = 78
$ws.Range("F48").Value = 4421
$ws.Range("F49").Value = 580

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 118
$ws.Range("F15").Value = 20
$ws.Range("F26").Value = 15

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7559
$ws.Range("F4").Value = 766

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 7559
$ws.Range("F5").Value = 766
$ws.Range("F6").Value = 14196
$ws.Range("F7").Value = 14299
$ws.Range("F9").Value = 1399
$ws.Range("F11").Value = 984
$ws.Range("F12").Value = 118
$ws.Range("F16").Value = 436
$ws.Range("F17").Value = 1198
$ws.Range("F18").Value = 1828
$ws.Range("F21").Value = 3318
$ws.Range("F22").Value = 311
$ws.Range("F23").Value = 2391
$ws.Range("F24").Value = 592
$ws.Range("F27").Value = 1788
$ws.Range("F31").Value = 1069
$ws.Range("F32").Value = 1389
$ws.Range("F33").Value = 100
$ws.Range("F34").Value = 4810
$ws.Range("F35").Value = 4845
$ws.Range("F36").Value = 302
$ws.Range("F38").Value = 672
$ws.Range("F40").Value = 3287
$ws.Range("F42").Value = 337
$ws.Range("F43").Value = 106
$ws.Range("F45").Value = 78
$ws.Range("F46").Value = 4421
$ws.Range("F47").Value = 580
$ws.Range("F49").Value = 15

$wb.Save()